$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 7488.467
$ws.Range("I6").Value = 909.7692
$ws.Range("J6").Value = 50250
$ws.Range("K6").Value = 2729.3076
$ws.Range("L6").Value = 150750
$ws.Range("M6").Value = -2617.3076
$ws.Range("N6").Value = -150974

$ws.Range("H33").Value = 299.32257
$ws.Range("I33").Value = 112.77778
$ws.Range("J33").Value = 1558.5
$ws.Range("K33").Value = 112.77778
$ws.Range("L33").Value = 1558.5
$ws.Range("M33").Value = 116.22222
$ws.Range("N33").Value = -2016.5

$ws.Range("H98").Value = 79109.2
$ws.Range("I98").Value = 1187.8572
$ws.Range("J98").Value = 260925.67
$ws.Range("K98").Value = 1187.8572
$ws.Range("L98").Value = 260925.67
$ws.Range("M98").Value = 310.1428000000001
$ws.Range("N98").Value = -263921.67

$ws.Range("H122").Value = 79109.2
$ws.Range("I122").Value = 1187.8572
$ws.Range("J122").Value = 260925.67
$ws.Range("K122").Value = 3563.5716
$ws.Range("L122").Value = 782777.01
$ws.Range("M122").Value = -1113.5716
$ws.Range("N122").Value = -787677.01

$ws.Range("H126").Value = 49800
$ws.Range("J126").Value = 49800
$ws.Range("L126").Value = 49800
$ws.Range("N126").Value = -59680

$ws.Range("H129").Value = 1058.9487
$ws.Range("J129").Value = 855.94116
$ws.Range("L129").Value = 2567.82348
$ws.Range("N129").Value = -12567.82348

$ws.Range("H137").Value = 2086357.1
$ws.Range("I137").Value = 9631622
$ws.Range("J137").Value = 4904.8276
$ws.Range("K137").Value = 28894866
$ws.Range("L137").Value = 14714.4828
$ws.Range("M137").Value = -28892316
$ws.Range("N137").Value = -19814.4828

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11849.075
$ws.Range("I32").Value = 12391.052
$ws.Range("J32").Value = 10339.286
$ws.Range("K32").Value = 12391.052
$ws.Range("L32").Value = 10339.286
$ws.Range("M32").Value = -12104.052
$ws.Range("N32").Value = -10913.286

$ws.Range("H48").Value = 243000
$ws.Range("J48").Value = 243000
$ws.Range("L48").Value = 243000
$ws.Range("N48").Value = -243768

$ws.Range("H122").Value = 2250
$ws.Range("I122").Value = 2142.8572
$ws.Range("K122").Value = 6428.571599999999
$ws.Range("M122").Value = -3978.571599999999

$ws.Range("H132").Value = 10419325
$ws.Range("I132").Value = 16131239
$ws.Range("J132").Value = 3482
$ws.Range("K132").Value = 48393717
$ws.Range("L132").Value = 10446
$ws.Range("M132").Value = -48391187
$ws.Range("N132").Value = -15506

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1179.84
$ws.Range("I94").Value = 1205.8889
$ws.Range("K94").Value = 1205.8889
$ws.Range("M94").Value = -754.8888999999999

$ws.Range("H107").Value = 3391.5557
$ws.Range("I107").Value = 2252.75
$ws.Range("J107").Value = 4302.6
$ws.Range("K107").Value = 2252.75
$ws.Range("L107").Value = 4302.6
$ws.Range("M107").Value = -332.75
$ws.Range("N107").Value = -8142.6

$ws.Range("H134").Value = 2795.2964
$ws.Range("I134").Value = 2611.5
$ws.Range("J134").Value = 4265.6665
$ws.Range("K134").Value = 7834.5
$ws.Range("L134").Value = 12796.9995
$ws.Range("M134").Value = -5299.5
$ws.Range("N134").Value = -17866.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4389697
$ws.Range("I31").Value = 1972.6364
$ws.Range("J31").Value = 6177288
$ws.Range("K31").Value = 1972.6364
$ws.Range("L31").Value = 6177288
$ws.Range("M31").Value = -1677.6364
$ws.Range("N31").Value = -6177878

$ws.Range("H34").Value = 4389697
$ws.Range("I34").Value = 1972.6364
$ws.Range("J34").Value = 6177288
$ws.Range("K34").Value = 1972.6364
$ws.Range("L34").Value = 6177288
$ws.Range("M34").Value = -1770.6364
$ws.Range("N34").Value = -6177692

$ws.Range("H58").Value = 2688.45
$ws.Range("I58").Value = 2725
$ws.Range("J58").Value = 2672.7856
$ws.Range("K58").Value = 2725
$ws.Range("L58").Value = 2672.7856
$ws.Range("M58").Value = -2522
$ws.Range("N58").Value = -3078.7856

$ws.Range("H105").Value = 1996.5186
$ws.Range("I105").Value = 2153.818
$ws.Range("J105").Value = 1304.4
$ws.Range("K105").Value = 2153.818
$ws.Range("L105").Value = 1304.4
$ws.Range("M105").Value = -406.8180000000002
$ws.Range("N105").Value = -4798.4

$ws.Range("H107").Value = 566.9091
$ws.Range("I107").Value = 343.6316
$ws.Range("J107").Value = 869.9286
$ws.Range("K107").Value = 343.6316
$ws.Range("L107").Value = 869.9286
$ws.Range("M107").Value = 1576.3684
$ws.Range("N107").Value = -4709.9286

$ws.Range("H136").Value = 2688.45
$ws.Range("I136").Value = 2725
$ws.Range("J136").Value = 2672.7856
$ws.Range("K136").Value = 8175
$ws.Range("L136").Value = 8018.3568
$ws.Range("M136").Value = -5625
$ws.Range("N136").Value = -13118.3568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 5500
$ws.Range("J25").Value = 5500
$ws.Range("L25").Value = 16500
$ws.Range("N25").Value = -16838

$ws.Range("H29").Value = 15284.571
$ws.Range("I29").Value = 1180
$ws.Range("J29").Value = 25863
$ws.Range("K29").Value = 3540
$ws.Range("L29").Value = 77589
$ws.Range("M29").Value = -3263
$ws.Range("N29").Value = -78143

$ws.Range("H30").Value = 5500
$ws.Range("J30").Value = 5500
$ws.Range("L30").Value = 16500
$ws.Range("N30").Value = -16704

$ws.Range("H35").Value = 1380.4
$ws.Range("J35").Value = 1380.4
$ws.Range("L35").Value = 4141.200000000001
$ws.Range("N35").Value = -4717.200000000001

$ws.Range("H36").Value = 400
$ws.Range("I36").Value = 400
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1200
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1031
$ws.Range("N36").ClearContents()

$ws.Range("H64").Value = 2857.1428
$ws.Range("I64").Value = 2500
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 7500
$ws.Range("L64").Value = 15000
$ws.Range("M64").Value = -7230
$ws.Range("N64").Value = -15540

$ws.Range("H67").Value = 2857.1428
$ws.Range("I67").Value = 2500
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 7500
$ws.Range("L67").Value = 15000
$ws.Range("M67").Value = -6564
$ws.Range("N67").Value = -16872

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1534.3636
$ws.Range("I122").Value = 1543
$ws.Range("J122").Value = 1519.25
$ws.Range("K122").Value = 4629
$ws.Range("L122").Value = 4557.75
$ws.Range("M122").Value = -2179
$ws.Range("N122").Value = -9457.75

$ws.Range("H132").Value = 29414544
$ws.Range("I132").Value = 43480256
$ws.Range("J132").Value = 4417.4546
$ws.Range("K132").Value = 130440768
$ws.Range("L132").Value = 13252.3638
$ws.Range("M132").Value = -130438238
$ws.Range("N132").Value = -18312.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 66670310
$ws.Range("I7").Value = 111112890
$ws.Range("J7").Value = 6451.6665
$ws.Range("K7").Value = 111112890
$ws.Range("L7").Value = 6451.6665
$ws.Range("M7").Value = -111112778
$ws.Range("N7").Value = -6675.6665

$ws.Range("H40").Value = 3528.3225
$ws.Range("I40").Value = 2976.7827
$ws.Range("J40").Value = 5114
$ws.Range("K40").Value = 2976.7827
$ws.Range("L40").Value = 5114
$ws.Range("M40").Value = -2840.7827
$ws.Range("N40").Value = -5386

$ws.Range("H126").Value = 66670310
$ws.Range("I126").Value = 111112890
$ws.Range("J126").Value = 6451.6665
$ws.Range("K126").Value = 333338670
$ws.Range("L126").Value = 19354.9995
$ws.Range("M126").Value = -333336200
$ws.Range("N126").Value = -24294.9995

$ws.Range("H132").Value = 3430.8635
$ws.Range("I132").Value = 2208.2727
$ws.Range("J132").Value = 4653.4546
$ws.Range("K132").Value = 6624.8181
$ws.Range("L132").Value = 13960.3638
$ws.Range("M132").Value = -4094.8181
$ws.Range("N132").Value = -19020.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 5014
$ws.Range("I21").Value = 5014
$ws.Range("K21").Value = 5014
$ws.Range("M21").Value = -4779

$ws.Range("H35").Value = 5014
$ws.Range("I35").Value = 5014
$ws.Range("K35").Value = 5014
$ws.Range("M35").Value = -4724

$ws.Range("H132").Value = 946886.25
$ws.Range("I132").Value = 1500461.1
$ws.Range("J132").Value = 2552.7058
$ws.Range("K132").Value = 4501383.300000001
$ws.Range("L132").Value = 7658.117400000001
$ws.Range("M132").Value = -4498853.300000001
$ws.Range("N132").Value = -12718.1174
